# Automatische test-sync: 2025-08-01 23:16:50
$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

# Append the new log row (row 5)
$logs.Range("A5").Value = "Kun jij dit even regelen?"
$logs.Range("B5").Value = "mailmind.test@zohomail.eu"
$logs.Range("C5").Value = "Testmail #1: Kun jij dit even regelen?"
$logs.Range("D5").Value = "Overig"
$logs.Range("E5").Value = "Beste klant,`nBedankt voor je e-mail. Om je vraag beter te kunnen beantwoorden, heb ik meer details nodig. Kunt u beschrijven waar u specifiek hulp bij nodig heeft? Als u meer informatie geeft, kan ik u beter van dienst zijn.`nMet vriendelijke groet,`n[Naam]`nE-mailassistent"
$logs.Range("F5").Value = "2025-08-01 23:16:44"
$logs.Range("G5").Value = "Ja"
$logs.Range("H5").Value = "Nee"
$logs.Range("I5").Value = "Ja"
$logs.Range("J5").Value = "Nee"

# Update the Dashboard summary count for "Overig"
$dashboard.Range("B2").Value = 4

# Extend the conditional formatting ranges to include the new row 5
$ranges = @("D2:D4", "G2:G4", "H2:H4", "I2:I4", "J2:J4")
foreach ($rangeAddr in $ranges) {
    $col = $rangeAddr.Substring(0, 1)
    $fcs = $logs.Range($rangeAddr).FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($logs.Range("$col`2:$col`5"))
    }
}

# Keep row 5's height in line with the other data rows (avoid the implicit
# autofit bump from the multi-line Antwoord cell so the row stays uniform)
$logs.Rows.Item(5).RowHeight = $logs.Rows.Item(4).RowHeight
